$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the Date property and shorten the Description text ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2026-01-16T13:49:34+00:00"
$wsMeta.Range("B12").Value = "Statut clinique du patient"

# --- Elements sheet: the base resource row ("Definition" column) shares the
#     same underlying text as the Metadata "Description" cell, so it picks
#     up the same corrected wording. ---
$wsElem = $wb.Worksheets.Item("Elements")
$wsElem.Range("M2").Value = "Statut clinique du patient"
